$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 additional period rows before the final (bottom-bordered) data row ---
# Old layout: rows 16-21 "middle" style, row 22 "bottom" style (last period), rows 27-28 footer
# New layout needs 9 period rows total (16-24), so insert 2 rows at 22:23
$ws.Rows("22:23").Insert()

# Copy formatting (styles/borders) from the template middle row (21) into the two new rows
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update worker identity ---
$ws.Range("B22:B23").Value = "CC"
$ws.Range("C16:C24").Value = "1068587640"
$ws.Range("D16:D24").Value = "LEANIS VILLALBA VASQUEZ"

# --- Update period codes (column E) for the 9 period rows ---
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2503"
$ws.Range("E18").Value = "2502"
$ws.Range("E19").Value = "2501"
$ws.Range("E20").Value = "2412"
$ws.Range("E21").Value = "2411"
$ws.Range("E22").Value = "2410"
$ws.Range("E23").Value = "2409"
$ws.Range("E24").Value = "2408"

# --- Update amounts (Valor Mora / Salario Basico) per period ---
$ws.Range("F16:F24").Value = 56760
$ws.Range("G16:G24").Value = 1419000

# --- Update summary fields ---
$ws.Range("E11").Value = 510840
$ws.Range("F13").Value = 9

# Column D ("Nombre Trabajador") is bestFit; refresh its width now that the new name is longer
# (closest achievable approximation of Excel's recalculated bestFit width in this environment)
$ws.Columns("D:D").ColumnWidth = 24.25

Write-Host "Done"
